$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the empty column C (previously unused, reserved for "kappa" values)
# so that the former column D ("V4-A3-S1_k") becomes column C and the
# former column E ("V4-A3-S1_m") becomes column D.
$ws.Columns("C").Delete()

# Reflect the resulting selection, matching how Excel leaves the selection
# on the column that was just deleted.
$ws.Range("C1").Activate()
$ws.Range("C1:C1048576").Select()
